$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Birds")
$ws2 = $wb.Worksheets.Item("Cages")

# --- Birds sheet: row 2 ---
$ws1.Range("A2").Value = 1232
$ws1.Range("B2").Value = "a"
$ws1.Range("C2").Value = "a"
$ws1.Range("D2").Value = 45061.766585648147
$ws1.Range("D2").NumberFormat = "m/d/yy h:mm"
$ws1.Range("E2").Value = "Female"
$ws1.Range("F2").Value = "d1a2s"
$ws1.Range("G2").Value = 12
$ws1.Range("H2").Value = 12
$ws1.Range("I2").Value = "s"
$ws1.Range("J2").Value = "s"
$ws1.Range("K2").Value = "s"

# --- Birds sheet: row 3 ---
$ws1.Range("A3").Value = 123
$ws1.Range("B3").Value = "a"
$ws1.Range("C3").Value = "a"
$ws1.Range("D3").Value = 45061.801041666666
$ws1.Range("D3").NumberFormat = "m/d/yy h:mm"
$ws1.Range("E3").Value = "Female"
$ws1.Range("F3").Value = "a1"
$ws1.Range("G3").Value = 1232
$ws1.Range("H3").Value = 1232
$ws1.Range("I3").Value = "a"
$ws1.Range("J3").Value = "a"
$ws1.Range("K3").Value = "a"

# --- Cages sheet: row 2 (new cage added, referenced by the bird in Birds!F6) ---
$ws2.Range("A2").Value = "adsa2"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = "Iron"

# --- Birds sheet: row 4 ---
$ws1.Range("A4").Value = 123111
$ws1.Range("B4").Value = "q"
$ws1.Range("C4").Value = "d"
$ws1.Range("D4").Value = 45061.869386574072
$ws1.Range("D4").NumberFormat = "m/d/yy h:mm"
$ws1.Range("E4").Value = "Male"
$ws1.Range("F4").Value = "12c"
$ws1.Range("G4").Value = 11
$ws1.Range("H4").Value = 1
$ws1.Range("I4").Value = "a"
$ws1.Range("J4").Value = "a"
$ws1.Range("K4").Value = "a"

# --- Birds sheet: row 5 ---
$ws1.Range("A5").Value = 123112
$ws1.Range("B5").Value = "das"
$ws1.Range("C5").Value = "das"
$ws1.Range("D5").Value = 45061.869756944441
$ws1.Range("D5").NumberFormat = "m/d/yy h:mm"
$ws1.Range("E5").Value = "Female"
$ws1.Range("F5").Value = "sda1"
$ws1.Range("G5").Value = 11
$ws1.Range("H5").Value = 111
$ws1.Range("I5").Value = "dsa"
$ws1.Range("J5").Value = "asd"
$ws1.Range("K5").Value = "das"

# --- Birds sheet: row 6 ---
$ws1.Range("A6").Value = 122
$ws1.Range("B6").Value = "das"
$ws1.Range("C6").Value = "das"
$ws1.Range("D6").Value = 45061.871377314812
$ws1.Range("D6").NumberFormat = "m/d/yy h:mm"
$ws1.Range("E6").Value = "Male"
$ws1.Range("F6").Value = "adsa2"
$ws1.Range("G6").Value = 12
$ws1.Range("H6").Value = 123
$ws1.Range("I6").Value = 1
$ws1.Range("J6").Value = 1
$ws1.Range("K6").Value = 1
